$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-20 12:47:53"

for ($r = 2; $r -le 24; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
